# Daily update at 8 AM UTC
# Appends the next day's row of win counts, moving the special
# "last row" date formatting from the old last row to the new one.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The previous last row (19) had the special "latest row" number format
# (style index 3 / numFmtId 167, "YYYY-MM-DD"). Since it's no longer the
# last row, give it the regular date format used by all other data rows
# (style index 2 / numFmtId 165, "YYYY-MM-DD HH:MM:SS") - same format as A2:A18.
$ws.Range("A19").NumberFormat = $ws.Range("A18").NumberFormat

# Add the new row of data for the next day.
$ws.Range("A20").Value = 45969
$ws.Range("B20").Value = 42
$ws.Range("C20").Value = 49
$ws.Range("D20").Value = 49

# The new last row gets the special "latest row" number format that A19
# used to have.
$ws.Range("A20").NumberFormat = "YYYY-MM-DD"
